# Update 18-Jan-2021, midday update.
# Applies the revised "Buku KAS HARIAN"/Sheet1 petty-cash entries: the
# earlier (pre-midday) entries for 11..15-Jan are replaced by a fresh,
# smaller set of entries dated 18-Jan, and the remainder of the day rows
# are cleared back out to blank (formulas in column E recompute on their
# own once the supporting cells are blank/changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Opening balance -------------------------------------------------
$ws.Range("E2").Value = 453525

# ---- Row 3: Wages Expense, now dated 18-Jan-2021 ----------------------
$ws.Range("A3").Value = 44214
$ws.Range("D3").Formula = "=45000"

# ---- Row 4: now "TAX - P.Tata", entered as a Credit -------------------
$ws.Range("B4").Value = "TAX - P.Tata"
$ws.Range("D4").Clear()
$ws.Range("C4").Formula = "=7300000+2800000"

# ---- Row 5: now "TRANSFER BCA" ----------------------------------------
$ws.Range("B5").Value = "TRANSFER BCA"
$ws.Range("D5").Formula = "=1787000+1149000"

# ---- Rows 6-46: all prior entries for this block are wiped out --------
# (the running-balance formulas in column E are untouched and simply
# recompute once their precedents disappear)
$ws.Range("A6:D46").Clear()

# ---- Restore the saved view state (freeze pane + active selection) ----
$ws.Range("D6").Select()
